$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 81682.336
$ws.Range("J21").Value = 49999
$ws.Range("L21").Value = 49999
$ws.Range("N21").Value = -50935

$ws.Range("H23").Value = 81682.336
$ws.Range("J23").Value = 49999
$ws.Range("L23").Value = 49999
$ws.Range("N23").Value = -50467

$ws.Range("H118").Value = 225.33333
$ws.Range("I118").Value = 245
$ws.Range("K118").Value = 735
$ws.Range("M118").Value = 922

$ws.Range("H138").Value = 2021.78
$ws.Range("I138").Value = 960.6279
$ws.Range("J138").Value = 2822.2983
$ws.Range("K138").Value = 2881.8837
$ws.Range("L138").Value = 8466.894899999999
$ws.Range("M138").Value = 2258.1163
$ws.Range("N138").Value = -18746.8949

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3347
$ws.Range("I2").Value = 1045
$ws.Range("K2").Value = 1045
$ws.Range("M2").Value = -932

$ws.Range("H32").Value = 12555.369
$ws.Range("I32").Value = 9012.135
$ws.Range("J32").Value = 17237.5
$ws.Range("K32").Value = 9012.135
$ws.Range("L32").Value = 17237.5
$ws.Range("M32").Value = -8725.135
$ws.Range("N32").Value = -17811.5

$ws.Range("H74").Value = 2491.7368
$ws.Range("I74").Value = 1788.6666
$ws.Range("K74").Value = 1788.6666
$ws.Range("M74").Value = -914.6666

$ws.Range("H77").Value = 2491.7368
$ws.Range("I77").Value = 1788.6666
$ws.Range("K77").Value = 8943.333000000001
$ws.Range("M77").Value = -4575.333000000001

$ws.Range("H116").Value = 3347
$ws.Range("I116").Value = 1045
$ws.Range("K116").Value = 1045
$ws.Range("M116").Value = 1249

$ws.Range("H132").Value = 1719.14
$ws.Range("I132").Value = 935.1905
$ws.Range("K132").Value = 2805.5715
$ws.Range("M132").Value = -275.5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3347
$ws.Range("I3").Value = 1045
$ws.Range("K3").Value = 1045
$ws.Range("M3").Value = -931

$ws.Range("H63").Value = 55000
$ws.Range("J63").Value = 55000
$ws.Range("L63").Value = 55000
$ws.Range("N63").Value = -56372

$ws.Range("H66").Value = 55000
$ws.Range("J66").Value = 55000
$ws.Range("L66").Value = 165000
$ws.Range("N66").Value = -171864

$ws.Range("H87").Value = 41800
$ws.Range("J87").Value = 41800
$ws.Range("L87").Value = 41800
$ws.Range("N87").Value = -44296

$ws.Range("H90").Value = 41800
$ws.Range("J90").Value = 41800
$ws.Range("L90").Value = 125400
$ws.Range("N90").Value = -137880

$ws.Range("H105").Value = 23811704
$ws.Range("J105").Value = 2677.75
$ws.Range("L105").Value = 2677.75
$ws.Range("N105").Value = -6171.75

$ws.Range("H107").Value = 1520
$ws.Range("I107").Value = 1396.75
$ws.Range("J107").Value = 2013
$ws.Range("K107").Value = 1396.75
$ws.Range("L107").Value = 2013
$ws.Range("M107").Value = 523.25
$ws.Range("N107").Value = -5853

$ws.Range("H134").Value = 2672.3845
$ws.Range("I134").Value = 1513.775
$ws.Range("J134").Value = 6534.4165
$ws.Range("K134").Value = 4541.325000000001
$ws.Range("L134").Value = 19603.2495
$ws.Range("M134").Value = -2006.325000000001
$ws.Range("N134").Value = -24673.2495

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2990.34
$ws.Range("I31").Value = 1081.7812
$ws.Range("J31").Value = 6383.3335
$ws.Range("K31").Value = 1081.7812
$ws.Range("L31").Value = 6383.3335
$ws.Range("M31").Value = -786.7811999999999
$ws.Range("N31").Value = -6973.3335

$ws.Range("H34").Value = 2990.34
$ws.Range("I34").Value = 1081.7812
$ws.Range("J34").Value = 6383.3335
$ws.Range("K34").Value = 1081.7812
$ws.Range("L34").Value = 6383.3335
$ws.Range("M34").Value = -879.7811999999999
$ws.Range("N34").Value = -6787.3335

$ws.Range("H141").Value = 15402.564
$ws.Range("J141").Value = 15402.564
$ws.Range("L141").Value = 15402.564
$ws.Range("N141").Value = -25762.564

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 769.5946
$ws.Range("I113").Value = 681.43475
$ws.Range("J113").Value = 914.4286
$ws.Range("K113").Value = 2044.30425
$ws.Range("L113").Value = 2743.2858
$ws.Range("M113").Value = 125.6957499999999
$ws.Range("N113").Value = -7083.2858

$ws.Range("H131").Value = 15625967
$ws.Range("I131").Value = 20000878
$ws.Range("J131").Value = 1284.5714
$ws.Range("K131").Value = 60002634
$ws.Range("L131").Value = 3853.7142
$ws.Range("M131").Value = -59997594
$ws.Range("N131").Value = -13933.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 16150
$ws.Range("J63").Value = 16150
$ws.Range("L63").Value = 16150
$ws.Range("N63").Value = -17522

$ws.Range("H66").Value = 16150
$ws.Range("J66").Value = 16150
$ws.Range("L66").Value = 48450
$ws.Range("N66").Value = -55314

$ws.Range("H70").Value = 6622.95
$ws.Range("I70").Value = 5738.25
$ws.Range("J70").Value = 7950
$ws.Range("K70").Value = 5738.25
$ws.Range("L70").Value = 7950
$ws.Range("M70").Value = -5468.25
$ws.Range("N70").Value = -8490

$ws.Range("H73").Value = 6622.95
$ws.Range("I73").Value = 5738.25
$ws.Range("J73").Value = 7950
$ws.Range("K73").Value = 5738.25
$ws.Range("L73").Value = 7950
$ws.Range("M73").Value = -4802.25
$ws.Range("N73").Value = -9822

$ws.Range("H80").Value = 27780376
$ws.Range("I80").Value = 250000000
$ws.Range("J80").Value = 2922.375
$ws.Range("K80").Value = 250000000
$ws.Range("L80").Value = 2922.375
$ws.Range("M80").Value = -249999002
$ws.Range("N80").Value = -4918.375

$ws.Range("H83").Value = 27780376
$ws.Range("I83").Value = 250000000
$ws.Range("J83").Value = 2922.375
$ws.Range("K83").Value = 1250000000
$ws.Range("L83").Value = 14611.875
$ws.Range("M83").Value = -1249995008
$ws.Range("N83").Value = -24595.875

$ws.Range("H102").Value = 2256.85
$ws.Range("I102").Value = 1965.5555
$ws.Range("J102").Value = 4878.5
$ws.Range("K102").Value = 1965.5555
$ws.Range("L102").Value = 4878.5
$ws.Range("M102").Value = -343.5554999999999
$ws.Range("N102").Value = -8122.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 53889.91
$ws.Range("J133").Value = 53889.91
$ws.Range("L133").Value = 53889.91
$ws.Range("N133").Value = -58949.91

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 39800
$ws.Range("J112").Value = 39800
$ws.Range("L112").Value = 39800
$ws.Range("N112").Value = -42754

$ws.Range("H132").Value = 7250343
$ws.Range("I132").Value = 4745.6
$ws.Range("J132").Value = 15876054
$ws.Range("K132").Value = 14236.8
$ws.Range("L132").Value = 47628162
$ws.Range("M132").Value = -11706.8
$ws.Range("N132").Value = -47633222
